# The authored change swaps the content of the two theme parts in the
# package (ppt/theme/theme1.xml <-> ppt/theme/theme2.xml):
#   - theme1.xml ("Office Theme" colours, only used by the Notes Master)
#     ends up holding what used to be the "Integral" theme.
#   - theme2.xml ("Integral" colours, used by the Slide Master / the
#     presentation's main design) ends up holding the "Office Theme"
#     colours that used to live in theme1.xml.
# The font scheme and format scheme (fills/lines/effects) are byte
# identical between the two theme parts already, so the only
# observable difference is the theme colour scheme (the 12 colour
# slots: dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# The design (colour scheme) that is reachable through the PowerPoint
# object model is the presentation's single Slide Master design, so
# we repaint its 12 theme colours to the "Office Theme" palette that
# used to live in theme1.xml.

$p = $ppt.ActivePresentation
$scheme = $p.SlideMaster.Theme.ThemeColorScheme

# MSO theme colour slot order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6,
# 11 hlink, 12 folHlink. Values are the target "Office Theme" palette,
# expressed as 0xBBGGRR integers (the VBA/COM RGB() encoding).
$targetColors = @{
    1  = 0        # dk1      000000
    2  = 16777215 # lt1      FFFFFF
    3  = 6968388  # dk2      44546A
    4  = 15132391 # lt2      E7E6E6
    5  = 13998939 # accent1  5B9BD5
    6  = 3243501  # accent2  ED7D31
    7  = 10855845 # accent3  A5A5A5
    8  = 49407    # accent4  FFC000
    9  = 12874308 # accent5  4472C4
    10 = 4697456  # accent6  70AD47
    11 = 12673797 # hlink    0563C1
    12 = 7491477  # folHlink 954F72
}

for ($i = 1; $i -le 12; $i++) {
    $scheme.Item($i).RGB = $targetColors[$i]
}
